# Generate Report for Handoff
# Rotate the localized-report identifiers from the old run
# (c31a51ad-570c-4506-ae90-a7ed0620939e, handoff zip
# eb7927f160b0f9ba61b9ea80df91c1d849d28cb7) to the new run
# (759239cf-2818-4f34-9c84-0fae4df38b1c, handoff zip
# 229120e0fd521357ec57b6b626152d250f46545a), and bump the
# generation/handoff timestamps that moved with it.

$wb = $excel.ActiveWorkbook

$oldGuid = "c31a51ad-570c-4506-ae90-a7ed0620939e"
$newGuid = "759239cf-2818-4f34-9c84-0fae4df38b1c"

$oldHash = "eb7927f160b0f9ba61b9ea80df91c1d849d28cb7"
$newHash = "229120e0fd521357ec57b6b626152d250f46545a"

# Hyperlinks in this workbook all still resolve to the original
# GitHub blob URL (the target commit didn't move) - only the
# visible/display text changes to the new file name.
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40703e7bc0e2822c82e3f6e6a65df264f7f06ef5/e2e/$oldGuid.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-19 15:03:55"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# ---- zh-cn sheet ----
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-19 15:03:50"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# ---- de-de sheet ----
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
# H2 (Latest Handback DateTime) is unchanged by this commit.

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
